# Update GRU model and results
# Adds a "GRU" model row to each of the four service blocks (CPU Usage table
# in A:E and the mirrored Memory Usage table in G:K), reusing the first
# previously-blank row in each 6-row block (rows 7, 13, 19, 25) and filling
# in its results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row (within each block) that becomes the new "GRU" row, and the row right
# above it (the existing "Random Forest Regressor" row) whose formatting we
# clone for the new row.
$gruRows  = @(7, 13, 19, 25)
$rfRows   = @(6, 12, 18, 24)

# CPU Usage (left table, columns C/D/E) and Memory Usage (right table,
# columns I/J/K) values for the new GRU row, per block.
$cpuValues = @{
    7  = @(0.9987, 0.9978, 0.9938)
    13 = @(0.9993, 0.9992, 0.9973)
    19 = @(0.9979, 0.9989, 0.9974)
    25 = @(0.9712, 0.9738, 0.9692)
}
$memValues = @{
    7  = @(0.4664, 0.4053, 0.3039)
    13 = @(0.5156, 0.4883, 0.4335)
    19 = @(0.6078, 0.5939, 0.4859)
    25 = @(0.9921, 0.9924, 0.985)
}

for ($i = 0; $i -lt $gruRows.Length; $i++) {
    $gru = $gruRows[$i]
    $rf  = $rfRows[$i]

    # Clone the formatting (font, border, alignment, number format, etc.)
    # of the "Random Forest Regressor" row onto the new "GRU" row, for both
    # the left (A:E) and right (G:K) tables.
    $ws.Range("A${rf}:E${rf}").Copy()
    $ws.Range("A${gru}:E${gru}").PasteSpecial(-4122)
    $ws.Range("G${rf}:K${rf}").Copy()
    $ws.Range("G${gru}:K${gru}").PasteSpecial(-4122)

    # Model name label.
    $ws.Range("B$gru").Value = "GRU"
    $ws.Range("H$gru").Value = "GRU"

    # CPU Usage results (left table).
    $cpu = $cpuValues[$gru]
    $ws.Range("C$gru").Value = $cpu[0]
    $ws.Range("D$gru").Value = $cpu[1]
    $ws.Range("E$gru").Value = $cpu[2]

    # Memory Usage results (right table).
    $mem = $memValues[$gru]
    $ws.Range("I$gru").Value = $mem[0]
    $ws.Range("J$gru").Value = $mem[1]
    $ws.Range("K$gru").Value = $mem[2]
}

$excel.CutCopyMode = 0

# Move the active selection, matching the author's final cursor position.
[void]$ws.Range("L22").Select()
